$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 11:59"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 7679908
$ws.Range("C4").Value = 264
$ws.Range("D4").Value = 4895291
$ws.Range("E4").Value = 2569578
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 215039

# Row 19: 'Banglades' -> 'Banglades'
$ws.Range("B19").Value = 371631
$ws.Range("C19").Value = 1499
$ws.Range("D19").Value = 284833
$ws.Range("E19").Value = 81393
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 5405

# Row 27: 'Israel' -> 'Israel'
$ws.Range("B27").Value = 274423
$ws.Range("C27").Value = 2114
$ws.Range("D27").Value = 208819
$ws.Range("E27").Value = 63833
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = 1771

# Row 41: 'Egipto' -> 'Polonia'
$ws.Range("A41").Value = "Polonia"
$ws.Range("B41").Value = 104316
$ws.Range("C41").Value = 2236
$ws.Range("D41").Value = 74158
$ws.Range("E41").Value = 27441
$ws.Range("G41").Value = 58
$ws.Range("H41").Value = 2717

# Row 42: 'Polonia' -> 'Egipto'
$ws.Range("A42").Value = "Egipto"
$ws.Range("B42").Value = 103781
$ws.Range("D42").Value = 97398
$ws.Range("E42").Value = 393
$ws.Range("H42").Value = 5990

# Row 43: 'Oman' -> 'Oman'
$ws.Range("B43").Value = 102648
$ws.Range("C43").Value = 834
$ws.Range("D43").Value = 91275
$ws.Range("E43").Value = 10383
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 990

# Row 62: 'Suiza' -> 'Suiza'
$ws.Range("D62").Value = 47300
$ws.Range("E62").Value = 6554

# Row 65: 'Austria' -> 'Austria'
$ws.Range("B65").Value = 49819
$ws.Range("C65").Value = 923
$ws.Range("D65").Value = 39790
$ws.Range("E65").Value = 9207
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 822

# Row 72: 'Kenia' -> 'Afganistan'
$ws.Range("A72").Value = "Afganistan"
$ws.Range("B72").Value = 39486
$ws.Range("C72").Value = 64
$ws.Range("D72").Value = 32977
$ws.Range("E72").Value = 5042
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 1467

# Row 73: 'Afganistan' -> 'Kenia'
$ws.Range("A73").Value = "Kenia"
$ws.Range("B73").Value = 39449
$ws.Range("D73").Value = 27035
$ws.Range("E73").Value = 11679
$ws.Range("H73").Value = 735

# Row 90: 'Croacia' -> 'Croacia'
$ws.Range("B90").Value = 18084
$ws.Range("C90").Value = 287
$ws.Range("D90").Value = 16192
$ws.Range("E90").Value = 1588
$ws.Range("G90").Value = 4
$ws.Range("H90").Value = 304

# Row 99: 'Malasia' -> 'Malasia'
$ws.Range("B99").Value = 13504
$ws.Range("C99").Value = 691
$ws.Range("D99").Value = 10427
$ws.Range("E99").Value = 2936
$ws.Range("G99").Value = 4
$ws.Range("H99").Value = 141

# Row 102: 'Guinea' -> 'Finlandia'
$ws.Range("A102").Value = "Finlandia"
$ws.Range("B102").Value = 10929
$ws.Range("C102").Value = 227
$ws.Range("D102").Value = 8100
$ws.Range("E102").Value = 2483
$ws.Range("H102").Value = 346

# Row 103: 'Consejo Danes para los Refugiados' -> 'Guinea'
$ws.Range("A103").Value = "Guinea"
$ws.Range("B103").Value = 10800
$ws.Range("D103").Value = 10161
$ws.Range("E103").Value = 572
$ws.Range("H103").Value = 67

# Row 104: 'Finlandia' -> 'Consejo Danes para los Refugiados'
$ws.Range("A104").Value = "Consejo Danes para los Refugiados"
$ws.Range("B104").Value = 10778
$ws.Range("D104").Value = 10239
$ws.Range("E104").Value = 265
$ws.Range("H104").Value = 274

# Row 117: 'Eslovenia' -> 'Eslovenia'
$ws.Range("B117").Value = 6764
$ws.Range("C117").Value = 189
$ws.Range("D117").Value = 4399
$ws.Range("E117").Value = 2206
$ws.Range("G117").Value = 3
$ws.Range("H117").Value = 159

# Row 127: 'Hong Kong' -> 'Hong Kong'
$ws.Range("B127").Value = 5133
$ws.Range("C127").Value = 8
$ws.Range("D127").Value = 4875
$ws.Range("E127").Value = 153

# Row 138: 'Mayotte' -> 'Sri Lanka'
$ws.Range("A138").Value = "Sri Lanka"
$ws.Range("B138").Value = 3979
$ws.Range("C138").Value = 466
$ws.Range("D138").Value = 3266
$ws.Range("E138").Value = 700
$ws.Range("H138").Value = 13

# Row 139: 'Somalia' -> 'Mayotte'
$ws.Range("A139").Value = "Mayotte"
$ws.Range("B139").Value = 3892
$ws.Range("D139").Value = 2964
$ws.Range("E139").Value = 886
$ws.Range("H139").Value = 42

# Row 140: 'Sri Lanka' -> 'Somalia'
$ws.Range("A140").Value = "Somalia"
$ws.Range("B140").Value = 3745
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 3010
$ws.Range("E140").Value = 636
$ws.Range("H140").Value = 99

# Row 174: 'Papua Nueva Guinea' -> 'Papua Nueva Guinea'
$ws.Range("B174").Value = 541
$ws.Range("C174").Value = 1
$ws.Range("E174").Value = 7

# Row 175: 'Taiwan' -> 'Taiwan'
$ws.Range("B175").Value = 521
$ws.Range("C175").Value = 3
$ws.Range("E175").Value = 29

# Row 194: 'Brunei' -> 'Brunei'
$ws.Range("D194").Value = 143
$ws.Range("E194").Value = 0

# Row 195: 'Liechtenstein' -> 'Liechtenstein'
$ws.Range("B195").Value = 130
$ws.Range("C195").Value = 3
$ws.Range("E195").Value = 13

# Row 215: 'Islas Malvinas' -> 'Montserrat'
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# Row 216: 'Montserrat' -> 'Islas Malvinas'
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
